$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamps (col A) and load forecast values (col B) for rows 2-93
$ws.Range("A2").Value = 45431
$ws.Range("B2").Value = 5000
$ws.Range("A3").Value = 45431.01041666666
$ws.Range("B3").Value = 4940
$ws.Range("A4").Value = 45431.02083333334
$ws.Range("B4").Value = 4880
$ws.Range("A5").Value = 45431.03125
$ws.Range("B5").Value = 4840
$ws.Range("A6").Value = 45431.04166666666
$ws.Range("B6").Value = 4800
$ws.Range("A7").Value = 45431.05208333334
$ws.Range("B7").Value = 4770
$ws.Range("A8").Value = 45431.0625
$ws.Range("B8").Value = 4740
$ws.Range("A9").Value = 45431.07291666666
$ws.Range("B9").Value = 4720
$ws.Range("A10").Value = 45431.08333333334
$ws.Range("B10").Value = 4700
$ws.Range("A11").Value = 45431.09375
$ws.Range("B11").Value = 4690
$ws.Range("A12").Value = 45431.10416666666
$ws.Range("B12").Value = 4680
$ws.Range("A13").Value = 45431.11458333334
$ws.Range("B13").Value = 4670
$ws.Range("A14").Value = 45431.125
$ws.Range("B14").Value = 4660
$ws.Range("A15").Value = 45431.13541666666
$ws.Range("B15").Value = 4650
$ws.Range("A16").Value = 45431.14583333334
$ws.Range("B16").Value = 4650
$ws.Range("A17").Value = 45431.15625
$ws.Range("B17").Value = 4650
$ws.Range("A18").Value = 45431.16666666666
$ws.Range("B18").Value = 4650
$ws.Range("A19").Value = 45431.17708333334
$ws.Range("B19").Value = 4650
$ws.Range("A20").Value = 45431.1875
$ws.Range("B20").Value = 4650
$ws.Range("A21").Value = 45431.19791666666
$ws.Range("B21").Value = 4640
$ws.Range("A22").Value = 45431.20833333334
$ws.Range("B22").Value = 4630
$ws.Range("A23").Value = 45431.21875
$ws.Range("B23").Value = 4620
$ws.Range("A24").Value = 45431.22916666666
$ws.Range("B24").Value = 4610
$ws.Range("A25").Value = 45431.23958333334
$ws.Range("B25").Value = 4600
$ws.Range("A26").Value = 45431.25
$ws.Range("B26").Value = 4600
$ws.Range("A27").Value = 45431.26041666666
$ws.Range("B27").Value = 4600
$ws.Range("A28").Value = 45431.27083333334
$ws.Range("B28").Value = 4600
$ws.Range("A29").Value = 45431.28125
$ws.Range("B29").Value = 4600
$ws.Range("A30").Value = 45431.29166666666
$ws.Range("B30").Value = 4600
$ws.Range("A31").Value = 45431.30208333334
$ws.Range("B31").Value = 4600
$ws.Range("A32").Value = 45431.3125
$ws.Range("B32").Value = 4600
$ws.Range("A33").Value = 45431.32291666666
$ws.Range("B33").Value = 4590
$ws.Range("A34").Value = 45431.33333333334
$ws.Range("B34").Value = 4560
$ws.Range("A35").Value = 45431.34375
$ws.Range("B35").Value = 4530
$ws.Range("A36").Value = 45431.35416666666
$ws.Range("B36").Value = 4480
$ws.Range("A37").Value = 45431.36458333334
$ws.Range("B37").Value = 4430
$ws.Range("A38").Value = 45431.375
$ws.Range("B38").Value = 4380
$ws.Range("A39").Value = 45431.38541666666
$ws.Range("B39").Value = 4320
$ws.Range("A40").Value = 45431.39583333334
$ws.Range("B40").Value = 4270
$ws.Range("A41").Value = 45431.40625
$ws.Range("B41").Value = 4230
$ws.Range("A42").Value = 45431.41666666666
$ws.Range("B42").Value = 4190
$ws.Range("A43").Value = 45431.42708333334
$ws.Range("B43").Value = 4150
$ws.Range("A44").Value = 45431.4375
$ws.Range("B44").Value = 4120
$ws.Range("A45").Value = 45431.44791666666
$ws.Range("B45").Value = 4100
$ws.Range("A46").Value = 45431.45833333334
$ws.Range("B46").Value = 4070
$ws.Range("A47").Value = 45431.46875
$ws.Range("B47").Value = 4050
$ws.Range("A48").Value = 45431.47916666666
$ws.Range("B48").Value = 4020
$ws.Range("A49").Value = 45431.48958333334
$ws.Range("B49").Value = 4000
$ws.Range("A50").Value = 45431.5
$ws.Range("B50").Value = 3970
$ws.Range("A51").Value = 45431.51041666666
$ws.Range("B51").Value = 3950
$ws.Range("A52").Value = 45431.52083333334
$ws.Range("B52").Value = 3940
$ws.Range("A53").Value = 45431.53125
$ws.Range("B53").Value = 3940
$ws.Range("A54").Value = 45431.54166666666
$ws.Range("B54").Value = 3940
$ws.Range("A55").Value = 45431.55208333334
$ws.Range("B55").Value = 3950
$ws.Range("A56").Value = 45431.5625
$ws.Range("B56").Value = 3970
$ws.Range("A57").Value = 45431.57291666666
$ws.Range("B57").Value = 3990
$ws.Range("A58").Value = 45431.58333333334
$ws.Range("B58").Value = 4010
$ws.Range("A59").Value = 45431.59375
$ws.Range("B59").Value = 4030
$ws.Range("A60").Value = 45431.60416666666
$ws.Range("B60").Value = 4060
$ws.Range("A61").Value = 45431.61458333334
$ws.Range("B61").Value = 4100
$ws.Range("A62").Value = 45431.625
$ws.Range("B62").Value = 4150
$ws.Range("A63").Value = 45431.63541666666
$ws.Range("B63").Value = 4210
$ws.Range("A64").Value = 45431.64583333334
$ws.Range("B64").Value = 4280
$ws.Range("A65").Value = 45431.65625
$ws.Range("B65").Value = 4370
$ws.Range("A66").Value = 45431.66666666666
$ws.Range("B66").Value = 4470
$ws.Range("A67").Value = 45431.67708333334
$ws.Range("B67").Value = 4560
$ws.Range("A68").Value = 45431.6875
$ws.Range("B68").Value = 4650
$ws.Range("A69").Value = 45431.69791666666
$ws.Range("B69").Value = 4730
$ws.Range("A70").Value = 45431.70833333334
$ws.Range("B70").Value = 4810
$ws.Range("A71").Value = 45431.71875
$ws.Range("B71").Value = 4890
$ws.Range("A72").Value = 45431.72916666666
$ws.Range("B72").Value = 4980
$ws.Range("A73").Value = 45431.73958333334
$ws.Range("B73").Value = 5100
$ws.Range("A74").Value = 45431.75
$ws.Range("B74").Value = 5220
$ws.Range("A75").Value = 45431.76041666666
$ws.Range("B75").Value = 5350
$ws.Range("A76").Value = 45431.77083333334
$ws.Range("B76").Value = 5470
$ws.Range("A77").Value = 45431.78125
$ws.Range("B77").Value = 5570
$ws.Range("A78").Value = 45431.79166666666
$ws.Range("B78").Value = 5660
$ws.Range("A79").Value = 45431.80208333334
$ws.Range("B79").Value = 5760
$ws.Range("A80").Value = 45431.8125
$ws.Range("B80").Value = 5850
$ws.Range("A81").Value = 45431.82291666666
$ws.Range("B81").Value = 5970
$ws.Range("A82").Value = 45431.83333333334
$ws.Range("B82").Value = 6030
$ws.Range("A83").Value = 45431.84375
$ws.Range("B83").Value = 6050
$ws.Range("A84").Value = 45431.85416666666
$ws.Range("B84").Value = 6050
$ws.Range("A85").Value = 45431.86458333334
$ws.Range("B85").Value = 6000
$ws.Range("A86").Value = 45431.875
$ws.Range("B86").Value = 5890
$ws.Range("A87").Value = 45431.88541666666
$ws.Range("B87").Value = 5830
$ws.Range("A88").Value = 45431.89583333334
$ws.Range("B88").Value = 5720
$ws.Range("A89").Value = 45431.90625
$ws.Range("B89").Value = 5560
$ws.Range("A90").Value = 45431.91666666666
$ws.Range("B90").Value = 5310
$ws.Range("A91").Value = 45431.92708333334
$ws.Range("B91").Value = 5180
$ws.Range("A92").Value = 45431.9375
$ws.Range("B92").Value = 5080
$ws.Range("A93").Value = 45431.94791666666
$ws.Range("B93").Value = 4970

# Remove the now-unused trailing rows 94-98 (one day of data -> one fewer day)
$ws.Rows("94:98").Delete()

